# Refresh the crypto price/volume columns (D = Price, E = Volume(1h))
# to the latest scraped snapshot, row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.115.72"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "'1.894.67"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'306.62"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.5209"
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("D8").Value = "'0.3759"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "'0.07259"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "'21.15"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").Value = "'0.8996"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "'0.08198"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "'1.967.55"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "'96.26"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "'5.302"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'0.000008598"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "'14.59"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'27.149.09"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "'5.081"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "'6.415"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").Value = "'2.319"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").Value = "'148.25"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").Value = "'18.18"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'1.735"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "'115.14"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'4.798"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "'4.856"
$ws.Range("E30").Value = "  -2.81%  "
$ws.Range("D31").Value = "'0.09203"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "'0.05023"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").Value = "'0.7924"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").Value = "'1.218"
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").Value = "'3.434"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("D36").Value = "'2.957"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").Value = "'2.613"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").Value = "'0.5714"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "'0.01994"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").Value = "'9.034"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "'6.552"
$ws.Range("E42").Value = "  -1.18%  "
$ws.Range("D43").Value = "'116.35"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").Value = "'0.1513"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "'0.4860"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'10.07"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'1.619"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").Value = "'63.57"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").Value = "'0.05925"
$ws.Range("E51").Value = "  -0.45%  "
